# Actualización automática desde WSL
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17, column A: tiny correction to the stored timestamp fraction
$ws.Cells.Item(17, 1).Value = 45877.70853273148

# New row 18: next hourly reading appended to the log
$ws.Cells.Item(18, 1).Value = 45877.75021111657
$ws.Cells.Item(18, 2).Value = 2025
$ws.Cells.Item(18, 3).Value = 32
$ws.Cells.Item(18, 4).Value = 17.28
$ws.Cells.Item(18, 5).Value = 81.16
$ws.Cells.Item(18, 6).Value = 12.3
$ws.Cells.Item(18, 7).Value = 5.58
$ws.Cells.Item(18, 8).Value = "E"
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = "18:00:18"

# Match the date/time style used by the rest of column A
$ws.Cells.Item(18, 1).NumberFormat = $ws.Cells.Item(17, 1).NumberFormat
